$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 206 - this shifts existing rows 206:252 down to 207:253
# and carries over row formatting (e.g. the date number format on column D).
$ws.Rows.Item(206).Insert()

# Populate the newly inserted row 206 with the new data record.
$ws.Range("A206").Value = 3
$ws.Range("B206").Value = "Femacal de La Calera"
$ws.Range("C206").Value = "Coquimbo"
$ws.Range("D206").Value = 44543
$ws.Range("E206").Value = 5
$ws.Range("F206").Value = 100114013
$ws.Range("G206").Value = "Zanahoria"
$ws.Range("H206").Value = "Sin especificar"
$ws.Range("I206").Value = "Primera"
$ws.Range("J206").Value = 460
$ws.Range("K206").Value = 6000
$ws.Range("L206").Value = 6500
$ws.Range("M206").Value = 6283
$ws.Range("N206").Value = "$/saco 20 kilos"
$ws.Range("O206").Value = "Provincia de Quillota"
$ws.Range("P206").Value = 314
$ws.Range("Q206").Value = 20
$ws.Range("R206").Value = "Hortaliza"
